$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newApprovers = "ann-acm@arkcase.org,sally-acm@arkcase.org,samuel-acm@arkcase.org"

$ws.Range("I20").Value = $newApprovers
$ws.Range("I22").Value = $newApprovers
$ws.Range("I23").Value = $newApprovers
